$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value2 = 20.408218
$ws.Cells.Item(2, 8).Value2 = 61.224654
$ws.Cells.Item(2, 9).Value2 = 0.1108535210972707
$ws.Cells.Item(2, 10).Value2 = 0.1108535210972707
$ws.Cells.Item(2, 13).Value2 = 2.113523666666667
$ws.Cells.Item(2, 14).Value2 = 6.340571000000001
$ws.Cells.Item(2, 15).Value2 = 0.2651220308693004
$ws.Cells.Item(2, 16).Value2 = 0.2651220308693004
$ws.Cells.Item(2, 17).Value2 = 43.13325173749268
$ws.Cells.Item(2, 18).Value2 = 388.1992656374341
$ws.Cells.Item(2, 19).Value2 = 0.02938971064232124
$ws.Cells.Item(2, 20).Value2 = 0.02938971064232124

# Row 3
$ws.Cells.Item(3, 7).Value2 = 20.408218
$ws.Cells.Item(3, 8).Value2 = 61.224654
$ws.Cells.Item(3, 9).Value2 = 0.1108535210972707
$ws.Cells.Item(3, 10).Value2 = 0.1108535210972707
$ws.Cells.Item(3, 15).Value2 = 0.2869289465860668
$ws.Cells.Item(3, 16).Value2 = 0.2869289465860668
$ws.Cells.Item(3, 17).Value2 = 46.68106397378801
$ws.Cells.Item(3, 18).Value2 = 420.129575764092
$ws.Cells.Item(3, 19).Value2 = 0.03180708403379621
$ws.Cells.Item(3, 20).Value2 = 0.03180708403379621

# Row 4
$ws.Cells.Item(4, 7).Value2 = 20.408218
$ws.Cells.Item(4, 8).Value2 = 61.224654
$ws.Cells.Item(4, 9).Value2 = 0.1108535210972707
$ws.Cells.Item(4, 10).Value2 = 0.1108535210972707
$ws.Cells.Item(4, 13).Value2 = 1.164746666666667
$ws.Cells.Item(4, 14).Value2 = 3.49424
$ws.Cells.Item(4, 15).Value2 = 0.1461067158059967
$ws.Cells.Item(4, 16).Value2 = 0.1461067158059966
$ws.Cells.Item(4, 17).Value2 = 23.77040388810667
$ws.Cells.Item(4, 18).Value2 = 213.93363499296
$ws.Cells.Item(4, 19).Value2 = 0.01619644390305298
$ws.Cells.Item(4, 20).Value2 = 0.01619644390305298

# Row 5
$ws.Cells.Item(5, 7).Value2 = 20.408218
$ws.Cells.Item(5, 8).Value2 = 61.224654
$ws.Cells.Item(5, 9).Value2 = 0.1108535210972707
$ws.Cells.Item(5, 10).Value2 = 0.1108535210972707
$ws.Cells.Item(5, 13).Value2 = 2.406253666666667
$ws.Cells.Item(5, 14).Value2 = 7.218761
$ws.Cells.Item(5, 15).Value2 = 0.3018423067386362
$ws.Cells.Item(5, 16).Value2 = 0.3018423067386362
$ws.Cells.Item(5, 17).Value2 = 49.10734939263267
$ws.Cells.Item(5, 18).Value2 = 441.966144533694
$ws.Cells.Item(5, 19).Value2 = 0.03346028251810026
$ws.Cells.Item(5, 20).Value2 = 0.03346028251810026

# Row 6
$ws.Cells.Item(6, 9).Value2 = 0.2566851044076959
$ws.Cells.Item(6, 10).Value2 = 0.256685104407696
$ws.Cells.Item(6, 13).Value2 = 2.113523666666667
$ws.Cells.Item(6, 14).Value2 = 6.340571000000001
$ws.Cells.Item(6, 15).Value2 = 0.2651220308693004
$ws.Cells.Item(6, 16).Value2 = 0.2651220308693004
$ws.Cells.Item(6, 17).Value2 = 99.87651376420133
$ws.Cells.Item(6, 18).Value2 = 898.8886238778119
$ws.Cells.Item(6, 19).Value2 = 0.06805287617446674
$ws.Cells.Item(6, 20).Value2 = 0.06805287617446676

# Row 7
$ws.Cells.Item(7, 9).Value2 = 0.2566851044076959
$ws.Cells.Item(7, 10).Value2 = 0.256685104407696
$ws.Cells.Item(7, 15).Value2 = 0.2869289465860668
$ws.Cells.Item(7, 16).Value2 = 0.2869289465860668
$ws.Cells.Item(7, 19).Value2 = 0.07365038661203475
$ws.Cells.Item(7, 20).Value2 = 0.07365038661203477

# Row 8
$ws.Cells.Item(8, 9).Value2 = 0.2566851044076959
$ws.Cells.Item(8, 10).Value2 = 0.256685104407696
$ws.Cells.Item(8, 13).Value2 = 1.164746666666667
$ws.Cells.Item(8, 14).Value2 = 3.49424
$ws.Cells.Item(8, 15).Value2 = 0.1461067158059967
$ws.Cells.Item(8, 16).Value2 = 0.1461067158059966
$ws.Cells.Item(8, 17).Value2 = 55.04117995925332
$ws.Cells.Item(8, 18).Value2 = 495.3706196332799
$ws.Cells.Item(8, 19).Value2 = 0.03750341760132781
$ws.Cells.Item(8, 20).Value2 = 0.03750341760132781

# Row 9
$ws.Cells.Item(9, 9).Value2 = 0.2566851044076959
$ws.Cells.Item(9, 10).Value2 = 0.256685104407696
$ws.Cells.Item(9, 13).Value2 = 2.406253666666667
$ws.Cells.Item(9, 14).Value2 = 7.218761
$ws.Cells.Item(9, 15).Value2 = 0.3018423067386362
$ws.Cells.Item(9, 16).Value2 = 0.3018423067386362
$ws.Cells.Item(9, 17).Value2 = 113.7097403967213
$ws.Cells.Item(9, 18).Value2 = 1023.387663570492
$ws.Cells.Item(9, 19).Value2 = 0.07747842401986661
$ws.Cells.Item(9, 20).Value2 = 0.07747842401986663

# Row 10
$ws.Cells.Item(10, 7).Value2 = 85.307233
$ws.Cells.Item(10, 8).Value2 = 255.921699
$ws.Cells.Item(10, 9).Value2 = 0.4633725077375833
$ws.Cells.Item(10, 10).Value2 = 0.4633725077375833
$ws.Cells.Item(10, 13).Value2 = 2.113523666666667
$ws.Cells.Item(10, 14).Value2 = 6.340571000000001
$ws.Cells.Item(10, 15).Value2 = 0.2651220308693004
$ws.Cells.Item(10, 16).Value2 = 0.2651220308693004
$ws.Cells.Item(10, 17).Value2 = 180.2988558833477
$ws.Cells.Item(10, 18).Value2 = 1622.689702950129
$ws.Cells.Item(10, 19).Value2 = 0.1228502603003887
$ws.Cells.Item(10, 20).Value2 = 0.1228502603003887

# Row 11
$ws.Cells.Item(11, 7).Value2 = 85.307233
$ws.Cells.Item(11, 8).Value2 = 255.921699
$ws.Cells.Item(11, 9).Value2 = 0.4633725077375833
$ws.Cells.Item(11, 10).Value2 = 0.4633725077375833
$ws.Cells.Item(11, 15).Value2 = 0.2869289465860668
$ws.Cells.Item(11, 16).Value2 = 0.2869289465860668
$ws.Cells.Item(11, 17).Value2 = 195.128864318278
$ws.Cells.Item(11, 18).Value2 = 1756.159778864502
$ws.Cells.Item(11, 19).Value2 = 0.1329549855220888
$ws.Cells.Item(11, 20).Value2 = 0.1329549855220888

# Row 12
$ws.Cells.Item(12, 7).Value2 = 85.307233
$ws.Cells.Item(12, 8).Value2 = 255.921699
$ws.Cells.Item(12, 9).Value2 = 0.4633725077375833
$ws.Cells.Item(12, 10).Value2 = 0.4633725077375833
$ws.Cells.Item(12, 13).Value2 = 1.164746666666667
$ws.Cells.Item(12, 14).Value2 = 3.49424
$ws.Cells.Item(12, 15).Value2 = 0.1461067158059967
$ws.Cells.Item(12, 16).Value2 = 0.1461067158059966
$ws.Cells.Item(12, 17).Value2 = 99.36131527930665
$ws.Cells.Item(12, 18).Value2 = 894.2518375137599
$ws.Cells.Item(12, 19).Value2 = 0.06770183530032707
$ws.Cells.Item(12, 20).Value2 = 0.06770183530032706

# Row 13
$ws.Cells.Item(13, 7).Value2 = 85.307233
$ws.Cells.Item(13, 8).Value2 = 255.921699
$ws.Cells.Item(13, 9).Value2 = 0.4633725077375833
$ws.Cells.Item(13, 10).Value2 = 0.4633725077375833
$ws.Cells.Item(13, 13).Value2 = 2.406253666666667
$ws.Cells.Item(13, 14).Value2 = 7.218761
$ws.Cells.Item(13, 15).Value2 = 0.3018423067386362
$ws.Cells.Item(13, 16).Value2 = 0.3018423067386362
$ws.Cells.Item(13, 17).Value2 = 205.2708421994377
$ws.Cells.Item(13, 18).Value2 = 1847.437579794939
$ws.Cells.Item(13, 19).Value2 = 0.1398654266147787
$ws.Cells.Item(13, 20).Value2 = 0.1398654266147787

# Row 14
$ws.Cells.Item(14, 7).Value2 = 31.12938966666666
$ws.Cells.Item(14, 8).Value2 = 93.38816899999999
$ws.Cells.Item(14, 9).Value2 = 0.16908886675745
$ws.Cells.Item(14, 10).Value2 = 0.16908886675745
$ws.Cells.Item(14, 13).Value2 = 2.113523666666667
$ws.Cells.Item(14, 14).Value2 = 6.340571000000001
$ws.Cells.Item(14, 15).Value2 = 0.2651220308693004
$ws.Cells.Item(14, 16).Value2 = 0.2651220308693004
$ws.Cells.Item(14, 17).Value2 = 65.79270178938879
$ws.Cells.Item(14, 18).Value2 = 592.134316104499
$ws.Cells.Item(14, 19).Value2 = 0.04482918375212368
$ws.Cells.Item(14, 20).Value2 = 0.04482918375212368

# Row 15
$ws.Cells.Item(15, 7).Value2 = 31.12938966666666
$ws.Cells.Item(15, 8).Value2 = 93.38816899999999
$ws.Cells.Item(15, 9).Value2 = 0.16908886675745
$ws.Cells.Item(15, 10).Value2 = 0.16908886675745
$ws.Cells.Item(15, 15).Value2 = 0.2869289465860668
$ws.Cells.Item(15, 16).Value2 = 0.2869289465860668
$ws.Cells.Item(15, 17).Value2 = 71.20430752428466
$ws.Cells.Item(15, 18).Value2 = 640.8387677185619
$ws.Cells.Item(15, 19).Value2 = 0.04851649041814694
$ws.Cells.Item(15, 20).Value2 = 0.04851649041814694

# Row 16
$ws.Cells.Item(16, 7).Value2 = 31.12938966666666
$ws.Cells.Item(16, 8).Value2 = 93.38816899999999
$ws.Cells.Item(16, 9).Value2 = 0.16908886675745
$ws.Cells.Item(16, 10).Value2 = 0.16908886675745
$ws.Cells.Item(16, 13).Value2 = 1.164746666666667
$ws.Cells.Item(16, 14).Value2 = 3.49424
$ws.Cells.Item(16, 15).Value2 = 0.1461067158059967
$ws.Cells.Item(16, 16).Value2 = 0.1461067158059966
$ws.Cells.Item(16, 17).Value2 = 36.25785284961778
$ws.Cells.Item(16, 18).Value2 = 326.3206756465599
$ws.Cells.Item(16, 19).Value2 = 0.02470501900128879
$ws.Cells.Item(16, 20).Value2 = 0.02470501900128878

# Row 17
$ws.Cells.Item(17, 7).Value2 = 31.12938966666666
$ws.Cells.Item(17, 8).Value2 = 93.38816899999999
$ws.Cells.Item(17, 9).Value2 = 0.16908886675745
$ws.Cells.Item(17, 10).Value2 = 0.16908886675745
$ws.Cells.Item(17, 13).Value2 = 2.406253666666667
$ws.Cells.Item(17, 14).Value2 = 7.218761
$ws.Cells.Item(17, 15).Value2 = 0.3018423067386362
$ws.Cells.Item(17, 16).Value2 = 0.3018423067386362
$ws.Cells.Item(17, 17).Value2 = 74.90520802651211
$ws.Cells.Item(17, 18).Value2 = 674.1468722386089
$ws.Cells.Item(17, 19).Value2 = 0.05103817358589062
$ws.Cells.Item(17, 20).Value2 = 0.05103817358589061

